$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "64.365.63"
Set-TextValue $ws.Range("E2") "  -2.50%  "
Set-TextValue $ws.Range("D3") "3.340.62"
Set-TextValue $ws.Range("E3") "  -4.20%  "
Set-TextValue $ws.Range("E4") "  +0.02%  "
Set-TextValue $ws.Range("D5") "552.62"
Set-TextValue $ws.Range("E5") "  -5.45%  "
Set-TextValue $ws.Range("D6") "175.57"
Set-TextValue $ws.Range("E6") "  -1.37%  "
Set-TextValue $ws.Range("D7") "0.618"
Set-TextValue $ws.Range("E7") "  -2.23%  "
Set-TextValue $ws.Range("B8") "LidoStakedEther"
Set-TextValue $ws.Range("C8") "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
Set-TextValue $ws.Range("D8") "3.330.91"
Set-TextValue $ws.Range("E8") "  -4.24%  "
Set-TextValue $ws.Range("B9") "USDC"
Set-TextValue $ws.Range("C9") "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
Set-TextValue $ws.Range("D9") "1.00"
Set-TextValue $ws.Range("E9") "  -0.01%  "
Set-TextValue $ws.Range("D10") "0.627"
Set-TextValue $ws.Range("E10") "  -1.43%  "
Set-TextValue $ws.Range("E11") "  +1.18%  "
Set-TextValue $ws.Range("D12") "54.72"
Set-TextValue $ws.Range("E12") "  -2.45%  "
Set-TextValue $ws.Range("D13") "0.0000272"
Set-TextValue $ws.Range("E13") "  -2.18%  "
Set-TextValue $ws.Range("D14") "9.05"
Set-TextValue $ws.Range("E14") "  -2.43%  "
Set-TextValue $ws.Range("D15") "3.870.60"
Set-TextValue $ws.Range("E15") "  -4.27%  "
Set-TextValue $ws.Range("D16") "18.31"
Set-TextValue $ws.Range("E16") "  +0.11%  "
Set-TextValue $ws.Range("D17") "0.118"
Set-TextValue $ws.Range("E17") "  -2.96%  "
Set-TextValue $ws.Range("D18") "3.334.37"
Set-TextValue $ws.Range("E18") "  -4.30%  "
Set-TextValue $ws.Range("B19") "Uniswap"
Set-TextValue $ws.Range("C19") "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws.Range("D19") "11.74"
Set-TextValue $ws.Range("E19") "  -2.68%  "
Set-TextValue $ws.Range("B20") "WrappedBTC"
Set-TextValue $ws.Range("C20") "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue $ws.Range("D20") "64.226.98"
Set-TextValue $ws.Range("E20") "  -2.63%  "
Set-TextValue $ws.Range("D21") "0.979"
Set-TextValue $ws.Range("E21") "  -3.14%  "
Set-TextValue $ws.Range("D22") "433.18"
Set-TextValue $ws.Range("E22") "  +4.86%  "
Set-TextValue $ws.Range("D23") "5.12"
Set-TextValue $ws.Range("E23") "  +16.77%  "
Set-TextValue $ws.Range("D24") "4.06"
Set-TextValue $ws.Range("E24") "  -5.92%  "
Set-TextValue $ws.Range("B25") "InternetComputer(DFINITY)"
Set-TextValue $ws.Range("C25") "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D25") "13.47"
Set-TextValue $ws.Range("E25") "  +0.20%  "
Set-TextValue $ws.Range("B26") "Litecoin"
Set-TextValue $ws.Range("C26") "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D26") "84.30"
Set-TextValue $ws.Range("E26") "  -0.52%  "
Set-TextValue $ws.Range("D27") "10.75"
Set-TextValue $ws.Range("E27") "  -3.01%  "
Set-TextValue $ws.Range("D28") "2.83"
Set-TextValue $ws.Range("E28") "  -1.25%  "
Set-TextValue $ws.Range("D29") "8.73"
Set-TextValue $ws.Range("E29") "  -5.30%  "
Set-TextValue $ws.Range("D30") "29.73"
Set-TextValue $ws.Range("E30") "  -1.91%  "
Set-TextValue $ws.Range("D31") "6.67"
Set-TextValue $ws.Range("E31") "  -0.21%  "
Set-TextValue $ws.Range("D32") "11.48"
Set-TextValue $ws.Range("E32") "  -2.31%  "
Set-TextValue $ws.Range("D33") "579.82"
Set-TextValue $ws.Range("E33") "  -2.68%  "
Set-TextValue $ws.Range("E34") "  -2.98%  "
Set-TextValue $ws.Range("D35") "58.20"
Set-TextValue $ws.Range("E35") "  -4.57%  "
Set-TextValue $ws.Range("E36") "  +0.00%  "
Set-TextValue $ws.Range("E37") "  -8.05%  "
Set-TextValue $ws.Range("D38") "3.49"
Set-TextValue $ws.Range("E38") "  -2.81%  "
Set-TextValue $ws.Range("D39") "35.66"
Set-TextValue $ws.Range("E39") "  -3.43%  "
Set-TextValue $ws.Range("D40") "0.0₃0752"
Set-TextValue $ws.Range("E40") "  -5.42%  "
Set-TextValue $ws.Range("D41") "0.367"
Set-TextValue $ws.Range("E41") "  -4.67%  "
Set-TextValue $ws.Range("D42") "3.109.10"
Set-TextValue $ws.Range("E42") "  -3.56%  "
Set-TextValue $ws.Range("D43") "0.998"
Set-TextValue $ws.Range("E43") "  -0.13%  "
Set-TextValue $ws.Range("E44") "  -6.34%  "
Set-TextValue $ws.Range("D45") "3.21"
Set-TextValue $ws.Range("E45") "  -3.65%  "
Set-TextValue $ws.Range("D46") "0.0408"
Set-TextValue $ws.Range("E46") "  -2.84%  "
Set-TextValue $ws.Range("D47") "2.46"
Set-TextValue $ws.Range("E47") "  -3.55%  "
Set-TextValue $ws.Range("E48") "  -2.12%  "
Set-TextValue $ws.Range("D49") "2.60"
Set-TextValue $ws.Range("E49") "  -2.33%  "
Set-TextValue $ws.Range("B50") "THORChain"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
Set-TextValue $ws.Range("D50") "8.27"
Set-TextValue $ws.Range("E50") "  -4.17%  "
Set-TextValue $ws.Range("B51") "Monero"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D51") "134.98"
Set-TextValue $ws.Range("E51") "  -3.56%  "
